# Updated cryptos list (price + 1h volume change) values.
# Equivalent to the GitHub Actions data refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.111.03"
$ws.Range("E2").Value = "  +0.21%  "
# Row 3
$ws.Range("D3").Value = "2.416.44"
$ws.Range("E3").Value = "  -0.01%  "
# Row 4
$ws.Range("E4").Value = "  -0.03%  "
# Row 5
$ws.Range("D5").Value = "'553.66"
$ws.Range("E5").Value = "  +0.28%  "
# Row 6
$ws.Range("D6").Value = "'136.79"
$ws.Range("E6").Value = "  -0.15%  "
# Row 7
$ws.Range("E7").Value = "  +0.00%  "
# Row 8
$ws.Range("E8").Value = "  +0.89%  "
# Row 9
$ws.Range("E9").Value = "  -0.73%  "
# Row 10
$ws.Range("D10").Value = "'5.68"
# Row 11
$ws.Range("E11").Value = "  -0.56%  "
# Row 12
$ws.Range("E12").Value = "  -1.12%  "
# Row 13
$ws.Range("D13").Value = "'24.78"
$ws.Range("E13").Value = "  -0.20%  "
# Row 14
$ws.Range("D14").Value = "2.850.49"
$ws.Range("E14").Value = "  +0.07%  "
# Row 15
$ws.Range("D15").Value = "60.015.92"
$ws.Range("E15").Value = "  +0.15%  "
# Row 17
$ws.Range("D17").Value = "2.422.30"
$ws.Range("E17").Value = "  +0.70%  "
# Row 18
$ws.Range("D18").Value = "'11.24"
# Row 19
$ws.Range("E19").Value = "  +2.62%  "
# Row 20
$ws.Range("D20").Value = "'327.36"
$ws.Range("E20").Value = "  -1.18%  "
# Row 21
$ws.Range("E21").Value = "  +1.09%  "
# Row 22
$ws.Range("E22").Value = "  +0.01%  "
# Row 23
$ws.Range("E23").Value = "  -1.09%  "
# Row 24
$ws.Range("E24").Value = "  +5.20%  "
# Row 25
$ws.Range("D25").Value = "'8.64"
$ws.Range("E25").Value = "  +0.44%  "
# Row 26
$ws.Range("E26").Value = "  +0.11%  "
# Row 27
$ws.Range("E27").Value = "  +5.98%  "
# Row 28
$ws.Range("E28").Value = "  -1.03%  "
# Row 29
$ws.Range("E29").Value = "  -0.18%  "
# Row 30
$ws.Range("D30").Value = "'170.60"
$ws.Range("E30").Value = "  +0.11%  "
# Row 31
$ws.Range("E31").Value = "  -1.66%  "
# Row 32
$ws.Range("E32").Value = "  +4.25%  "
# Row 33
$ws.Range("D33").Value = "'0.401"
$ws.Range("E33").Value = "  -3.05%  "
# Row 34
$ws.Range("D34").Value = "'18.51"
$ws.Range("E34").Value = "  -0.47%  "
# Row 36
$ws.Range("E36").Value = "  +2.78%  "
# Row 37
$ws.Range("D37").Value = "'4.24"
$ws.Range("E37").Value = "  +1.23%  "
# Row 38
$ws.Range("E38").Value = "  +0.07%  "
# Row 39
$ws.Range("D39").Value = "'322.21"
$ws.Range("E39").Value = "  +3.18%  "
# Row 40
$ws.Range("E40").Value = "  -0.96%  "
# Row 41
$ws.Range("D41").Value = "'145.86"
$ws.Range("E41").Value = "  +4.44%  "
# Row 42
$ws.Range("E42").Value = "  -0.77%  "
# Row 43
$ws.Range("D43").Value = "'0.0964"
# Row 44
$ws.Range("D44").Value = "'19.80"
$ws.Range("E44").Value = "  +2.67%  "
# Row 46
$ws.Range("E46").Value = "  +0.69%  "
# Row 47
$ws.Range("D47").Value = "'0.0222"
$ws.Range("E47").Value = "  -1.29%  "
# Row 48
$ws.Range("E48").Value = "  -0.11%  "
# Row 49
$ws.Range("D49").Value = "'1.58"
$ws.Range("E49").Value = "  -0.85%  "
# Row 50
$ws.Range("E50").Value = "  -0.68%  "
# Row 51
$ws.Range("D51").Value = "'0.938"
$ws.Range("E51").Value = "  -1.54%  "
